# Auto-generated Excel COM-interop script to apply diff changes (cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '24.878.84'
$ws.Range("E2").Value = '  +1.10%  '
$ws.Range("D3").Value = '1.708.17'
$ws.Range("E3").Value = '  +0.97%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.36'
$ws.Range("E5").Value = '  +0.60%  '
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4020'
$ws.Range("E7").Value = '  +3.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4045'
$ws.Range("E8").Value = '  +0.64%  '
$ws.Range("B9").Value = 'Polygon'
$ws.Range("C9").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.474'
$ws.Range("E9").Value = '  -2.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '53.66'
$ws.Range("E10").Value = '  +1.63%  '
$ws.Range("B11").Value = 'BinanceUSD'
$ws.Range("C11").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9997'
$ws.Range("E11").Value = '  -0.41%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08804'
$ws.Range("E12").Value = '  +1.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '26.31'
$ws.Range("E13").Value = '  +6.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.511'
$ws.Range("E14").Value = '  -1.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.003'
$ws.Range("E15").Value = '  +0.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001343'
$ws.Range("E16").Value = '  +0.43%  '
$ws.Range("D17").Value = '1.667.94'
$ws.Range("E17").Value = '  -1.17%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.51'
$ws.Range("E18").Value = '  -2.66%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07184'
$ws.Range("E19").Value = '  +1.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.99'
$ws.Range("E20").Value = '  +6.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.293'
$ws.Range("E21").Value = '  +0.67%  '
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("E23").Value = '  +1.88%  '
$ws.Range("D24").Value = '24.869.93'
$ws.Range("E24").Value = '  +1.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.343'
$ws.Range("E25").Value = '  -0.17%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.892'
$ws.Range("E26").Value = '  -3.74%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.473'
$ws.Range("E27").Value = '  +23.94%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.11'
$ws.Range("E28").Value = '  +2.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '161.50'
$ws.Range("E29").Value = '  +0.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '144.09'
$ws.Range("E30").Value = '  +5.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.310'
$ws.Range("E31").Value = '  -2.41%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.284'
$ws.Range("E32").Value = '  +14.70%  '
$ws.Range("B33").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C33").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D33").Value = '1.845.01'
$ws.Range("E33").Value = '  -0.42%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08684'
$ws.Range("E34").Value = '  -0.93%  '
$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.03191'
$ws.Range("E35").Value = '  +10.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.233'
$ws.Range("E36").Value = '  -2.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.030'
$ws.Range("E37").Value = '  -0.27%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2863'
$ws.Range("E38").Value = '  +5.60%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8417'
$ws.Range("E39").Value = '  +8.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '10.84'
$ws.Range("E40").Value = '  +0.98%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.09440'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.31'
$ws.Range("E42").Value = '  +1.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.480'
$ws.Range("E43").Value = '  +1.90%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.51'
$ws.Range("E44").Value = '  +5.41%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.739'
$ws.Range("E45").Value = '  +6.65%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7428'
$ws.Range("E46").Value = '  +3.95%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.221'
$ws.Range("E47").Value = '  +0.53%  '
$ws.Range("E48").Value = '  +1.83%  '
$ws.Range("E49").Value = '  -0.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '140.66'
$ws.Range("E50").Value = '  +2.08%  '
$ws.Range("E51").Value = '  +5.15%  '
